# refactor: ahora cada registro de excel es una instancia del modelo
#
# Adds a new column C ("STRING") to every row of the mapping sheet, marking
# each mapped field with its type. The handful of "section header" rows in
# column A (previously styled with the blue/underlined hyperlink-ish font)
# lose that styling in favour of plain text, since the distinctive styling
# now lives on the new column C cells instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 185

# Rows whose column-A cell currently carries the "accent" font (Cambria,
# underlined, blue) that should move off of column A.
$accentRows = @(2, 4, 6, 8, 21, 40, 55, 59, 85, 116, 118, 144, 149)

for ($row = 1; $row -le $lastRow; $row++) {
    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value = "STRING"

    # New column C reuses the accent font (Cambria 11, single underline, blue)
    $cFont = $cCell.Font
    $cFont.Name = "Cambria"
    $cFont.Size = 11
    $cFont.Underline = 2
    $cFont.Color = 16711680
    $cFont.Bold = $false
    $cFont.Italic = $false
}

foreach ($row in $accentRows) {
    $aFont = $ws.Cells.Item($row, 1).Font
    $aFont.Name = "Arial"
    $aFont.Size = 10
    $aFont.Underline = -4142
    $aFont.Color = 0
    $aFont.Bold = $false
    $aFont.Italic = $false
}

# Match the author's final cursor position / selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 169
$win.ScrollColumn = 1
$ws.Range("D188").Select()
